# Fix Training Data Issue (#48)
#
# The "Date" column was mislabeled for every team row because of the way NBA
# stats were shown on stats.nba.com (the scrape picked up data that was one
# calendar day off). Correct the stored date string from "6-24-2012-13" to
# the proper "2013-06-24" for every data row in the sheet.
#
# NOTE: Excel auto-recognizes strings shaped like "2013-06-24" as dates, so a
# plain Value assignment would silently convert the text into a date serial
# number instead of keeping it as literal text (which is how the column is
# stored, matching every other string column on the sheet). To avoid that
# implicit conversion we briefly force the cell to Text format while writing
# the corrected value, then clear the formatting again so the cell is left
# exactly as it was before (default/general formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "6-24-2012-13"
$newDate = "2013-06-24"

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1
$firstCol = $usedRange.Column
$lastCol = $firstCol + $usedRange.Columns.Count - 1

# Locate the "Date" column on the header row (row 1).
$dateCol = 0
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $header = $ws.Cells.Item($firstRow, $col).Text
    if ($header -eq "Date") {
        $dateCol = $col
        break
    }
}

if ($dateCol -eq 0) {
    # Fallback to the known column (BF) if the header could not be located.
    $dateCol = 58
}

for ($row = $firstRow + 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Text -eq $oldDate) {
        $cell.NumberFormat = "@"
        $cell.Value = $newDate
        $cell.ClearFormats()
    }
}
